$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 846, shifting rows 846-877 down to 848-879
$ws.Rows("846:847").Insert()

# Row 846: new weekly entry - "Primera" quality, 12-unit box
$ws.Cells.Item(846, 1).Value  = 3
$ws.Cells.Item(846, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(846, 3).Value  = "Coquimbo"
$ws.Cells.Item(846, 4).Value  = 44939
$ws.Cells.Item(846, 5).Value  = 5
$ws.Cells.Item(846, 6).Value  = "Fruta"
$ws.Cells.Item(846, 7).Value  = 100108
$ws.Cells.Item(846, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(846, 9).Value  = 100108005
$ws.Cells.Item(846, 10).Value = "Piña"
$ws.Cells.Item(846, 11).Value = "Caramelo"
$ws.Cells.Item(846, 12).Value = "Primera"
$ws.Cells.Item(846, 13).Value = 108
$ws.Cells.Item(846, 14).Value = 21000
$ws.Cells.Item(846, 15).Value = 21000
$ws.Cells.Item(846, 16).Value = 21000
$ws.Cells.Item(846, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(846, 18).Value = "Ecuador"
$ws.Cells.Item(846, 19).Value = 1750
$ws.Cells.Item(846, 20).Value = 12

# Row 847: new weekly entry - "Segunda" quality, 14-unit box
$ws.Cells.Item(847, 1).Value  = 3
$ws.Cells.Item(847, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(847, 3).Value  = "Coquimbo"
$ws.Cells.Item(847, 4).Value  = 44939
$ws.Cells.Item(847, 5).Value  = 5
$ws.Cells.Item(847, 6).Value  = "Fruta"
$ws.Cells.Item(847, 7).Value  = 100108
$ws.Cells.Item(847, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(847, 9).Value  = 100108005
$ws.Cells.Item(847, 10).Value = "Piña"
$ws.Cells.Item(847, 11).Value = "Caramelo"
$ws.Cells.Item(847, 12).Value = "Segunda"
$ws.Cells.Item(847, 13).Value = 108
$ws.Cells.Item(847, 14).Value = 21000
$ws.Cells.Item(847, 15).Value = 21000
$ws.Cells.Item(847, 16).Value = 21000
$ws.Cells.Item(847, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(847, 18).Value = "Ecuador"
$ws.Cells.Item(847, 19).Value = 1500
$ws.Cells.Item(847, 20).Value = 14
